$wb = $excel.ActiveWorkbook

# 1. On the "About" sheet, remove the explanatory paragraph about why two
#    versions of ETLE exist (rows 12-17), which also removes the now-unused
#    shared strings and shifts the "For more on this..." / URL rows up.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows("12:17").Delete()

# 2. On the "ETLE-capacity" sheet, change the exponent value in B2 from -1 to -5.
$wsCapacity = $wb.Worksheets.Item("ETLE-capacity")
$wsCapacity.Range("B2").Value = -5
